$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1807.138
$ws.Range("I43").Value = 1878.8422
$ws.Range("J43").Value = 1670.9
$ws.Range("K43").Value = 1878.8422
$ws.Range("L43").Value = 1670.9
$ws.Range("M43").Value = -1809.8422
$ws.Range("N43").Value = -1808.9
$ws.Range("H58").Value = 5006.467
$ws.Range("I58").Value = 664.25
$ws.Range("K58").Value = 1992.75
$ws.Range("M58").Value = -1842.75
$ws.Range("H69").Value = 14626.143
$ws.Range("I69").Value = 0.0
$ws.Range("K69").Value = 0.0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 14626.143
$ws.Range("I72").Value = 0.0
$ws.Range("K72").Value = 0.0
$ws.Range("M72").ClearContents()
$ws.Range("H74").Value = 13704.889
$ws.Range("I74").Value = 12293.0
$ws.Range("J74").Value = 25000.0
$ws.Range("K74").Value = 12293.0
$ws.Range("L74").Value = 25000.0
$ws.Range("M74").Value = -11357.0
$ws.Range("N74").Value = -26872.0
$ws.Range("H77").Value = 13704.889
$ws.Range("I77").Value = 12293.0
$ws.Range("J77").Value = 25000.0
$ws.Range("K77").Value = 61465.0
$ws.Range("L77").Value = 125000.0
$ws.Range("M77").Value = -56785.0
$ws.Range("N77").Value = -134360.0
$ws.Range("H98").Value = 1308.875
$ws.Range("I98").Value = 1062.8
$ws.Range("K98").Value = 1062.8
$ws.Range("M98").Value = 435.2
$ws.Range("H112").Value = 1760.0869
$ws.Range("I112").Value = 1100.0
$ws.Range("J112").Value = 1790.091
$ws.Range("K112").Value = 3300.0
$ws.Range("L112").Value = 5370.272999999999
$ws.Range("M112").Value = -2192.0
$ws.Range("N112").Value = -7586.272999999999
$ws.Range("H116").Value = 10518.182
$ws.Range("I116").Value = 5899.5
$ws.Range("K116").Value = 5899.5
$ws.Range("M116").Value = -2457.5
$ws.Range("H122").Value = 1308.875
$ws.Range("I122").Value = 1062.8
$ws.Range("K122").Value = 3188.4
$ws.Range("M122").Value = -738.3999999999996
$ws.Range("H132").Value = 2630.3142
$ws.Range("I132").Value = 2627.8064
$ws.Range("K132").Value = 7883.4192
$ws.Range("M132").Value = -5353.4192

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1140.2632
$ws.Range("I2").Value = 789.93335
$ws.Range("J2").Value = 2454.0
$ws.Range("K2").Value = 789.93335
$ws.Range("L2").Value = 2454.0
$ws.Range("M2").Value = -676.93335
$ws.Range("N2").Value = -2680.0
$ws.Range("H45").Value = 2056.05
$ws.Range("I45").Value = 1544.8
$ws.Range("K45").Value = 1544.8
$ws.Range("M45").Value = -1167.8
$ws.Range("H102").Value = 1893.85
$ws.Range("I102").Value = 1787.0588
$ws.Range("K102").Value = 1787.0588
$ws.Range("M102").Value = -165.0588
$ws.Range("H116").Value = 1140.2632
$ws.Range("I116").Value = 789.93335
$ws.Range("J116").Value = 2454.0
$ws.Range("K116").Value = 789.93335
$ws.Range("L116").Value = 2454.0
$ws.Range("M116").Value = 1504.06665
$ws.Range("N116").Value = -7042.0
$ws.Range("H122").Value = 4127.8276
$ws.Range("I122").Value = 1627.1111
$ws.Range("J122").Value = 5253.15
$ws.Range("K122").Value = 4881.3333
$ws.Range("L122").Value = 15759.45
$ws.Range("M122").Value = -2431.3333
$ws.Range("N122").Value = -20659.45
$ws.Range("H132").Value = 3062.9868
$ws.Range("I132").Value = 2200.6428
$ws.Range("J132").Value = 5477.55
$ws.Range("K132").Value = 6601.928400000001
$ws.Range("L132").Value = 16432.65
$ws.Range("M132").Value = -4071.928400000001
$ws.Range("N132").Value = -21492.65

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1140.2632
$ws.Range("I3").Value = 789.93335
$ws.Range("J3").Value = 2454.0
$ws.Range("K3").Value = 789.93335
$ws.Range("L3").Value = 2454.0
$ws.Range("M3").Value = -675.93335
$ws.Range("N3").Value = -2682.0
$ws.Range("H80").Value = 167.4
$ws.Range("I80").Value = 120.25
$ws.Range("J80").Value = 198.83333
$ws.Range("K80").Value = 120.25
$ws.Range("L80").Value = 198.83333
$ws.Range("M80").Value = 877.75
$ws.Range("N80").Value = -2194.83333
$ws.Range("H83").Value = 167.4
$ws.Range("I83").Value = 120.25
$ws.Range("J83").Value = 198.83333
$ws.Range("K83").Value = 601.25
$ws.Range("L83").Value = 994.1666499999999
$ws.Range("M83").Value = 4390.75
$ws.Range("N83").Value = -10978.16665
$ws.Range("H94").Value = 1279.9556
$ws.Range("I94").Value = 1324.6842
$ws.Range("K94").Value = 1324.6842
$ws.Range("M94").Value = -873.6841999999999
$ws.Range("H134").Value = 20828.158
$ws.Range("I134").Value = 2815.6667
$ws.Range("K134").Value = 8447.000100000001
$ws.Range("M134").Value = -5912.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4239.5356
$ws.Range("I16").Value = 3799.3
$ws.Range("J16").Value = 5340.125
$ws.Range("K16").Value = 3799.3
$ws.Range("L16").Value = 5340.125
$ws.Range("M16").Value = -3512.3
$ws.Range("N16").Value = -5914.125
$ws.Range("H31").Value = 41917.72
$ws.Range("I31").Value = 1184.0
$ws.Range("K31").Value = 1184.0
$ws.Range("M31").Value = -889.0
$ws.Range("H34").Value = 41917.72
$ws.Range("I34").Value = 1184.0
$ws.Range("K34").Value = 1184.0
$ws.Range("M34").Value = -982.0
$ws.Range("H63").Value = 25555.555
$ws.Range("J63").Value = 25555.555
$ws.Range("L63").Value = 25555.555
$ws.Range("N63").Value = -26927.555
$ws.Range("H66").Value = 25555.555
$ws.Range("J66").Value = 25555.555
$ws.Range("L66").Value = 76666.66500000001
$ws.Range("N66").Value = -83530.66500000001
$ws.Range("H94").Value = 1740.0
$ws.Range("J94").Value = 2266.6667
$ws.Range("L94").Value = 2266.6667
$ws.Range("N94").Value = -3168.6667
$ws.Range("H97").Value = 64180.332
$ws.Range("J97").Value = 64180.332
$ws.Range("L97").Value = 64180.332
$ws.Range("N97").Value = -66162.332
$ws.Range("H113").Value = 4239.5356
$ws.Range("I113").Value = 3799.3
$ws.Range("J113").Value = 5340.125
$ws.Range("K113").Value = 3799.3
$ws.Range("L113").Value = 5340.125
$ws.Range("M113").Value = -1629.3
$ws.Range("N113").Value = -9680.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 636155.9
$ws.Range("I5").Value = 53989.535
$ws.Range("K5").Value = 161968.605
$ws.Range("M5").Value = -161856.605
$ws.Range("H135").Value = 636155.9
$ws.Range("I135").Value = 53989.535
$ws.Range("K135").Value = 485905.8150000001
$ws.Range("M135").Value = -483370.8150000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0.0
$ws.Range("J53").Value = 0.0
$ws.Range("L53").Value = 0.0
$ws.Range("N53").ClearContents()
$ws.Range("H122").Value = 3486.2
$ws.Range("I122").Value = 2495.1333
$ws.Range("J122").Value = 6459.4
$ws.Range("K122").Value = 7485.3999
$ws.Range("L122").Value = 19378.2
$ws.Range("M122").Value = -5035.3999
$ws.Range("N122").Value = -24278.2
$ws.Range("H136").Value = 49123.082
$ws.Range("J136").Value = 49123.082
$ws.Range("L136").Value = 147369.246
$ws.Range("N136").Value = -152469.246

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 0.0
$ws.Range("J57").Value = 0.0
$ws.Range("L57").Value = 0.0
$ws.Range("N57").ClearContents()
$ws.Range("H136").Value = 181275.11
$ws.Range("I136").Value = 317251.6
$ws.Range("J136").Value = 7225.2
$ws.Range("K136").Value = 951754.7999999999
$ws.Range("L136").Value = 21675.6
$ws.Range("M136").Value = -949204.7999999999
$ws.Range("N136").Value = -26775.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 18279.545
$ws.Range("I4").Value = 22308.334
$ws.Range("K4").Value = 22308.334
$ws.Range("M4").Value = -22195.334
$ws.Range("H96").Value = 92433.18
$ws.Range("I96").Value = 144193.86
$ws.Range("J96").Value = 1852.0
$ws.Range("K96").Value = 144193.86
$ws.Range("L96").Value = 1852.0
$ws.Range("M96").Value = -142820.86
$ws.Range("N96").Value = -4598.0
$ws.Range("H119").Value = 99900.0
$ws.Range("J119").Value = 99900.0
$ws.Range("L119").Value = 99900.0
$ws.Range("N119").Value = -109576.0
$ws.Range("H122").Value = 26319724.0
$ws.Range("I122").Value = 35718144.0
$ws.Range("J122").Value = 4146.8
$ws.Range("K122").Value = 107154432.0
$ws.Range("L122").Value = 12440.4
$ws.Range("M122").Value = -107151982.0
$ws.Range("N122").Value = -17340.4
$ws.Range("H136").Value = 236377.2
$ws.Range("I136").Value = 228466.34
$ws.Range("K136").Value = 685399.02
$ws.Range("M136").Value = -685399.02
